$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "66.412.82"
$ws.Cells.Item(2, 5).Value = "  +0.00%  "

$ws.Cells.Item(3, 4).Value = "3.215.56"
$ws.Cells.Item(3, 5).Value = "  +0.90%  "

$ws.Cells.Item(4, 4).Value = "'0.999"
$ws.Cells.Item(4, 4).Style = "Normal"
$ws.Cells.Item(4, 5).Value = "  -0.08%  "

$ws.Cells.Item(5, 4).Value = "'609.12"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  +2.20%  "

$ws.Cells.Item(6, 4).Value = "'157.39"
$ws.Cells.Item(6, 4).Style = "Normal"

$ws.Cells.Item(7, 5).Value = "  +0.06%  "

$ws.Cells.Item(8, 4).Value = "3.214.60"
$ws.Cells.Item(8, 5).Value = "  +1.19%  "

$ws.Cells.Item(9, 4).Value = "'0.554"
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).Value = "  +0.51%  "

$ws.Cells.Item(10, 5).Value = "  +0.69%  "

$ws.Cells.Item(11, 4).Value = "'5.70"
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  -4.19%  "

$ws.Cells.Item(12, 5).Value = "  -2.56%  "

$ws.Cells.Item(13, 5).Value = "  +0.91%  "

$ws.Cells.Item(14, 4).Value = "'38.70"
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = "  -1.39%  "

$ws.Cells.Item(15, 4).Value = "3.743.01"
$ws.Cells.Item(15, 5).Value = "  +0.98%  "

$ws.Cells.Item(16, 4).Value = "66.481.84"
$ws.Cells.Item(16, 5).Value = "  +0.17%  "

$ws.Cells.Item(17, 4).Value = "'7.38"
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 5).Value = "  -1.38%  "

$ws.Cells.Item(18, 4).Value = "3.224.66"

$ws.Cells.Item(19, 4).Value = "'0.114"
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  +1.47%  "

$ws.Cells.Item(20, 4).Value = "'512.13"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  -0.82%  "

$ws.Cells.Item(21, 4).Value = "'15.23"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  -1.00%  "

$ws.Cells.Item(22, 4).Value = "'0.733"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  -0.70%  "

$ws.Cells.Item(23, 5).Value = "  -0.83%  "

$ws.Cells.Item(24, 4).Value = "'14.66"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  -1.93%  "

$ws.Cells.Item(25, 4).Value = "'85.15"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  -0.58%  "

$ws.Cells.Item(26, 5).Value = "  -0.03%  "

$ws.Cells.Item(27, 5).Value = "  +0.22%  "

$ws.Cells.Item(28, 4).Value = "'9.11"
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Cells.Item(28, 5).Value = "  -1.97%  "

$ws.Cells.Item(29, 4).Value = "'2.37"
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  +1.56%  "

$ws.Cells.Item(30, 5).Value = "  +39.55%  "

$ws.Cells.Item(31, 5).Value = "  +0.21%  "

$ws.Cells.Item(32, 5).Value = "  -1.79%  "

$ws.Cells.Item(33, 4).Value = "'28.26"
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  -0.29%  "

$ws.Cells.Item(34, 5).Value = "  +0.05%  "

$ws.Cells.Item(35, 5).Value = "  -4.65%  "

$ws.Cells.Item(36, 4).Value = "'6.52"
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  -0.06%  "

$ws.Cells.Item(37, 4).Value = "'503.77"
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  -0.72%  "

$ws.Cells.Item(38, 5).Value = "  +0.72%  "

$ws.Cells.Item(39, 4).Value = "0.0₃0779"
$ws.Cells.Item(39, 5).Value = "  +15.51%  "

$ws.Cells.Item(40, 5).Value = "  +6.72%  "

$ws.Cells.Item(41, 5).Value = "  -0.68%  "

$ws.Cells.Item(42, 4).Value = "'0.130"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  +1.80%  "

$ws.Cells.Item(43, 4).Value = "'8.75"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  -1.86%  "

$ws.Cells.Item(44, 4).Value = "'0.299"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  -1.34%  "

$ws.Cells.Item(45, 4).Value = "'2.45"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  +0.48%  "

$ws.Cells.Item(46, 4).Value = "2.919.89"
$ws.Cells.Item(46, 5).Value = "  +0.35%  "

$ws.Cells.Item(47, 4).Value = "'28.27"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  -1.27%  "

$ws.Cells.Item(48, 5).Value = "  +3.54%  "

$ws.Cells.Item(50, 4).Value = "'0.117"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  -0.76%  "

$ws.Cells.Item(51, 4).Value = "'122.10"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  -0.88%  "
